$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD (Wins), AE (Losses), AF (Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from the existing
# last header cell (AC1) onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season-record values (Wins/Losses/Ties) for every data row (2-41)
$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("AD$r").Value = 98
    $ws.Range("AE$r").Value = 64
    $ws.Range("AF$r").Value = 0
}

Write-Output "Done"
